$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 6
    3  = 2
    4  = -2
    5  = -2
    6  = 6
    7  = 7
    8  = -3
    9  = 1
    10 = -1
    11 = 11
    12 = -3
    14 = 1
    15 = 1
    16 = -2
    17 = -1
    20 = 2
    21 = -7
    22 = 1
    24 = -2
    25 = -3
    27 = 3
    28 = 4
    29 = 4
    30 = 2
    31 = 4
    32 = -1
    33 = -2
    34 = -4
    35 = -4
    37 = -2
    38 = 3
    39 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
